# df_tabela.xlsx -> relatorio de medias das notas
# Original sheet had 6 columns (aluno_id, created_at, nome, updated_at,
# curso_id, matricula_id). The report keeps only two columns: nome and the
# average grade per student. Delete the unneeded columns from the right so
# indices stay stable, leaving column A (ids) and column B (names) — then
# delete the dates column (B) so the remaining two columns are ids/names,
# and finally overwrite the values in place. This reuses the existing
# header style (bold/border/center) instead of minting a new one, and
# leaves the data rows with the sheet's default (unstyled) format, exactly
# like the original data rows did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(6).Delete()   # matricula_id
$ws.Columns.Item(5).Delete()   # curso_id
$ws.Columns.Item(4).Delete()   # updated_at
$ws.Columns.Item(2).Delete()   # created_at

# Remaining columns: A = aluno_id (ids), B = nome (names)
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "Média das Notas"

$ws.Range("A2").Value = "Big Smoke"
$ws.Range("B2").Value = 6.083333333333333

$ws.Range("A3").Value = "Cj"
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = "Cleber"
$ws.Range("B4").Value = 1.166666666666667

$ws.Range("A5").Value = "Jorge"
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = "Lamar"
$ws.Range("B6").Value = 2.5

# Column widths: COM's ColumnWidth adds the standard ~0.8333 padding unit
# relative to the raw stored width, so subtract it back out to land on the
# exact target widths (10 / 19).
$ws.Columns.Item(1).ColumnWidth = 10 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 19 - 0.8333333333333334
